# Literature review workbook: rename the original notes sheet to "KG ",
# then add two new sheets ("Paper 1 - 4 " and "Paper 5 - 9") that hold a
# compact per-paper summary table (Title / Main task / Architecture / Key
# Concepts / Training Data Set / Performance Measured on Data Set /
# Performance Metric / Performance Value).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the sole existing worksheet to "KG "
# ---------------------------------------------------------------------
$kg = $wb.Worksheets.Item(1)
$kg.Name = "KG "

# ---------------------------------------------------------------------
# 2. Add "Paper 1 - 4 " right after "KG "
# ---------------------------------------------------------------------
$p14 = $wb.Worksheets.Add($null, $kg)
$p14.Name = "Paper 1 - 4 "

$p14.Range("A1").Value = "Title"
$p14.Range("B1").Value = "Attention is all you need"
$p14.Range("C1").Value = "T"
$p14.Range("D1").Value = "SpanBert"
$p14.Range("E1").Value = "Personalized"

$p14.Range("A2").Value = "Main task"

$p14.Range("A3").Value = "Architecture"
$p14.Range("B3").Value = "Encoder  (6 layer Decoder Stack + <Some"
$p14.Range("E3").Value = "Bi LSTM"

$p14.Range("A4").Value = "Key Concepts"
$p14.Range("B4").Value = "Attention, Positional Encoding"

$p14.Range("A5").Value = "Training Data Set"
$p14.Range("B5").Value = "English German (4.5M Sentence Pairs, 37K tokens/English French (36M Sentences,25K tokens)"

$p14.Range("A6").Value = "Performance Measured on Data Set"
$p14.Range("B6").Value = "SQuaAD 1.1, SQuAD 2.0, MRQA"

$p14.Range("A7").Value = "Performance Metric"
$p14.Range("B7").Value = "Exact Match, F1 Score"

$p14.Range("A8").Value = "Performance Value"
$p14.Range("B8").Value = "1.1 - (EM 88.8, F1 94.6)`n2.0 - (EM -85.7, 88.7)`nMRQA - Avg (F1 -81.5))"

# Column B holds the long-form answers, so it gets wrapped text; column A
# (labels) and the scratch columns C:E stay plain. Row 2 has no answer yet,
# so only touch the cells that actually hold text.
$p14.Range("B1").WrapText = $true
$p14.Range("B3").WrapText = $true
$p14.Range("B4").WrapText = $true
$p14.Range("B5").WrapText = $true
$p14.Range("B6").WrapText = $true
$p14.Range("B7").WrapText = $true
$p14.Range("B8").WrapText = $true

$p14.Columns.Item(1).ColumnWidth = 27.93
$p14.Columns.Item(2).ColumnWidth = 27.83

$p14.Rows.Item(1).RowHeight = 16
$p14.Rows.Item(3).RowHeight = 32
$p14.Rows.Item(4).RowHeight = 16
$p14.Rows.Item(5).RowHeight = 48
$p14.Rows.Item(6).RowHeight = 16
$p14.Rows.Item(7).RowHeight = 16
$p14.Rows.Item(8).RowHeight = 48

$p14.Range("A1:A8").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. Add "Paper 5 - 9" right after "Paper 1 - 4 " (just the label column,
#    not yet filled in with per-paper data)
# ---------------------------------------------------------------------
$p59 = $wb.Worksheets.Add($null, $p14)
$p59.Name = "Paper 5 - 9"

$p59.Range("A1").Value = "Title"
$p59.Range("A2").Value = "Main task"
$p59.Range("A3").Value = "Architecture"
$p59.Range("A4").Value = "Key Concepts"
$p59.Range("A5").Value = "Training Data Set"
$p59.Range("A6").Value = "Performance Measured on Data Set"
$p59.Range("A7").Value = "Performance Metric"
$p59.Range("A8").Value = "Performance Value"

$p59.Columns.Item(1).ColumnWidth = 32.93

$p59.Range("A1:A8").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. Restore the selection on "KG " (it used to be on C7, now on the
#    merged detail cell B8:D8) and make "Paper 1 - 4 " the active tab,
#    matching where the author left off editing.
# ---------------------------------------------------------------------
$kg.Range("B8:D8").Select() | Out-Null
$p14.Activate()
